$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 34.973489275110438
$ws.Range("A2").Value = 40.638899258189845
$ws.Range("A3").Value = 41.69715590037071
$ws.Range("A4").Value = 42.143529676800085
$ws.Range("A5").Value = 26.257512368767443
$ws.Range("A6").Value = 25.639203030627129
$ws.Range("A7").Value = 38.997475802023288
$ws.Range("A8").Value = 36.670347194065116
$ws.Range("A9").Value = 34.230495268587916
$ws.Range("A10").Value = 37.698789617894214
$ws.Range("A11").Value = 26.148514567465352
